$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.06577466666666666
$ws.Range("H2").Value = 0.197324
$ws.Range("I2").Value = 0.1530524945763335
$ws.Range("J2").Value = 0.1530524945763335
$ws.Range("M2").Value = 0.02648366666666667
$ws.Range("N2").Value = 0.07945099999999999
$ws.Range("O2").Value = 0.001430039273477916
$ws.Range("P2").Value = 0.001430039273477917
$ws.Range("Q2").Value = 0.001741954347111111
$ws.Range("R2").Value = 0.015677589124
$ws.Range("S2").Value = 0.0002188710781479227
$ws.Range("T2").Value = 0.0002188710781479227
$ws.Range("G3").Value = 0.06577466666666666
$ws.Range("H3").Value = 0.197324
$ws.Range("I3").Value = 0.1530524945763335
$ws.Range("J3").Value = 0.1530524945763335
$ws.Range("O3").Value = 0.7016741634339546
$ws.Range("P3").Value = 0.7016741634339547
$ws.Range("Q3").Value = 0.8547208331395554
$ws.Range("R3").Value = 7.692487498256
$ws.Range("S3").Value = 0.1073929810933287
$ws.Range("T3").Value = 0.1073929810933287
$ws.Range("G4").Value = 0.06577466666666666
$ws.Range("H4").Value = 0.197324
$ws.Range("I4").Value = 0.1530524945763335
$ws.Range("J4").Value = 0.1530524945763335
$ws.Range("O4").Value = 0.2968957972925674
$ws.Range("P4").Value = 0.2968957972925675
$ws.Range("Q4").Value = 0.361653651284
$ws.Range("R4").Value = 3.254882861556
$ws.Range("S4").Value = 0.04544064240485689
$ws.Range("T4").Value = 0.04544064240485689
$ws.Range("H5").Value = 0.8943449999999999
$ws.Range("I5").Value = 0.6936902417438882
$ws.Range("J5").Value = 0.693690241743888
$ws.Range("M5").Value = 0.02648366666666667
$ws.Range("N5").Value = 0.07945099999999999
$ws.Range("O5").Value = 0.001430039273477916
$ws.Range("P5").Value = 0.001430039273477917
$ws.Range("Q5").Value = 0.007895178288333332
$ws.Range("R5").Value = 0.07105660459499999
$ws.Range("S5").Value = 0.0009920042893221501
$ws.Range("T5").Value = 0.0009920042893221501
$ws.Range("H6").Value = 0.8943449999999999
$ws.Range("I6").Value = 0.6936902417438882
$ws.Range("J6").Value = 0.693690241743888
$ws.Range("O6").Value = 0.7016741634339546
$ws.Range("P6").Value = 0.7016741634339547
$ws.Range("Q6").Value = 3.873909425686666
$ws.Range("S6").Value = 0.4867445200579404
$ws.Range("T6").Value = 0.4867445200579404
$ws.Range("H7").Value = 0.8943449999999999
$ws.Range("I7").Value = 0.6936902417438882
$ws.Range("J7").Value = 0.693690241743888
$ws.Range("O7").Value = 0.2968957972925674
$ws.Range("P7").Value = 0.2968957972925675
$ws.Range("S7").Value = 0.2059537173966255
$ws.Range("T7").Value = 0.2059537173966255
$ws.Range("I8").Value = 0.1532572636797784
$ws.Range("J8").Value = 0.1532572636797783
$ws.Range("M8").Value = 0.02648366666666667
$ws.Range("N8").Value = 0.07945099999999999
$ws.Range("O8").Value = 0.001430039273477916
$ws.Range("P8").Value = 0.001430039273477917
$ws.Range("Q8").Value = 0.001744284909777778
$ws.Range("R8").Value = 0.015698564188
$ws.Range("S8").Value = 0.0002191639060078437
$ws.Range("T8").Value = 0.0002191639060078437
$ws.Range("I9").Value = 0.1532572636797784
$ws.Range("J9").Value = 0.1532572636797783
$ws.Range("O9").Value = 0.7016741634339546
$ws.Range("P9").Value = 0.7016741634339547
$ws.Range("S9").Value = 0.1075366622826855
$ws.Range("T9").Value = 0.1075366622826855
$ws.Range("I10").Value = 0.1532572636797784
$ws.Range("J10").Value = 0.1532572636797783
$ws.Range("O10").Value = 0.2968957972925674
$ws.Range("P10").Value = 0.2968957972925675
$ws.Range("S10").Value = 0.04550143749108504
$ws.Range("T10").Value = 0.04550143749108504
